# Auto-generated Excel COM-interop script
# Applies numeric updates to the H-N "price/profit" columns across several
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) per the scheduled-runner sync.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2870.7925
$ws.Range("I74").Value = 2583.6365
$ws.Range("J74").Value = 3074.5806
$ws.Range("K74").Value = 2583.6365
$ws.Range("L74").Value = 3074.5806
$ws.Range("M74").Value = -1647.6365
$ws.Range("N74").Value = -4946.580599999999

$ws.Range("H77").Value = 2870.7925
$ws.Range("I77").Value = 2583.6365
$ws.Range("J77").Value = 3074.5806
$ws.Range("K77").Value = 12918.1825
$ws.Range("L77").Value = 15372.903
$ws.Range("M77").Value = -8238.182500000001
$ws.Range("N77").Value = -24732.903

$ws.Range("H138").Value = 2598.194
$ws.Range("I138").Value = 1675.8837
$ws.Range("J138").Value = 4250.6665
$ws.Range("K138").Value = 5027.6511
$ws.Range("L138").Value = 12751.9995
$ws.Range("M138").Value = 112.3489
$ws.Range("N138").Value = -23031.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1925.4814
$ws.Range("I88").Value = 2081.9092
$ws.Range("J88").Value = 1817.9375
$ws.Range("K88").Value = 2081.9092
$ws.Range("L88").Value = 1817.9375
$ws.Range("M88").Value = -1675.9092
$ws.Range("N88").Value = -2629.9375

$ws.Range("H91").Value = 1925.4814
$ws.Range("I91").Value = 2081.9092
$ws.Range("J91").Value = 1817.9375
$ws.Range("K91").Value = 2081.9092
$ws.Range("L91").Value = 1817.9375
$ws.Range("M91").Value = -677.9092000000001
$ws.Range("N91").Value = -4625.9375

$ws.Range("H96").Value = 25344
$ws.Range("J96").Value = 25344
$ws.Range("L96").Value = 25344
$ws.Range("N96").Value = -30836

$ws.Range("H134").Value = 32071.75
$ws.Range("J134").Value = 32071.75
$ws.Range("L134").Value = 32071.75
$ws.Range("N134").Value = -42211.75

$ws.Range("H139").Value = 85810
$ws.Range("J139").Value = 85810
$ws.Range("L139").Value = 85810
$ws.Range("N139").Value = -96090

$ws.Range("H140").Value = 44429
$ws.Range("J140").Value = 44429
$ws.Range("L140").Value = 44429
$ws.Range("N140").Value = -54789

$ws.Range("H141").Value = 40821.75
$ws.Range("J141").Value = 40821.75
$ws.Range("L141").Value = 40821.75
$ws.Range("N141").Value = -51181.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H134").Value = 4324.051
$ws.Range("I134").Value = 3664.3262
$ws.Range("K134").Value = 10992.9786
$ws.Range("M134").Value = -8457.9786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2494.375
$ws.Range("I62").Value = 2435
$ws.Range("J62").Value = 2530
$ws.Range("K62").Value = 2435
$ws.Range("L62").Value = 2530
$ws.Range("M62").Value = -1811
$ws.Range("N62").Value = -3778

$ws.Range("H65").Value = 2494.375
$ws.Range("I65").Value = 2435
$ws.Range("J65").Value = 2530
$ws.Range("K65").Value = 12175
$ws.Range("L65").Value = 12650
$ws.Range("M65").Value = -9055
$ws.Range("N65").Value = -18890

$ws.Range("H99").Value = 5646.6665
$ws.Range("I99").Value = 8633.333000000001
$ws.Range("J99").Value = 2660
$ws.Range("K99").Value = 8633.333000000001
$ws.Range("L99").Value = 2660
$ws.Range("M99").Value = -7135.333000000001
$ws.Range("N99").Value = -5656

$ws.Range("H126").Value = 5646.6665
$ws.Range("I126").Value = 8633.333000000001
$ws.Range("J126").Value = 2660
$ws.Range("K126").Value = 25899.999
$ws.Range("L126").Value = 7980
$ws.Range("M126").Value = -23429.999
$ws.Range("N126").Value = -12920

$ws.Range("H140").Value = 20000
$ws.Range("J140").Value = 20000
$ws.Range("L140").Value = 20000
$ws.Range("N140").Value = -30360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 836.6667
$ws.Range("I5").Value = 425.375
$ws.Range("J5").Value = 1223.7646
$ws.Range("K5").Value = 1276.125
$ws.Range("L5").Value = 3671.2938
$ws.Range("M5").Value = -1164.125
$ws.Range("N5").Value = -3895.2938

$ws.Range("H113").Value = 1894451
$ws.Range("I113").Value = 3030851.8
$ws.Range("J113").Value = 450
$ws.Range("K113").Value = 9092555.399999999
$ws.Range("L113").Value = 1350
$ws.Range("M113").Value = -9090385.399999999
$ws.Range("N113").Value = -5690

$ws.Range("H131").Value = 973.87
$ws.Range("I131").Value = 520
$ws.Range("J131").Value = 987.9072
$ws.Range("K131").Value = 1560
$ws.Range("L131").Value = 2963.7216
$ws.Range("M131").Value = 3480
$ws.Range("N131").Value = -13043.7216

$ws.Range("H135").Value = 836.6667
$ws.Range("I135").Value = 425.375
$ws.Range("J135").Value = 1223.7646
$ws.Range("K135").Value = 3828.375
$ws.Range("L135").Value = 11013.8814
$ws.Range("M135").Value = -1293.375
$ws.Range("N135").Value = -16083.8814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 13017.5
$ws.Range("J5").Value = 13017.5
$ws.Range("L5").Value = 13017.5
$ws.Range("N5").Value = -13241.5

$ws.Range("H133").Value = 62000
$ws.Range("J133").Value = 62000
$ws.Range("L133").Value = 62000
$ws.Range("N133").Value = -72120

$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10302000
$ws.Range("I2").Value = 12862500
$ws.Range("J2").Value = 60000
$ws.Range("K2").Value = 12862500
$ws.Range("L2").Value = 60000
$ws.Range("M2").Value = -12862388
$ws.Range("N2").Value = -60224

$ws.Range("H134").Value = 38966.668
$ws.Range("J134").Value = 38966.668
$ws.Range("L134").Value = 38966.668
$ws.Range("N134").Value = -49106.668

$ws.Range("H135").Value = 33201.707
$ws.Range("J135").Value = 33201.707
$ws.Range("L135").Value = 33201.707
$ws.Range("N135").Value = -43341.707

$ws.Range("H137").Value = 58830.77
$ws.Range("I137").Value = 64600
$ws.Range("J137").Value = 57100
$ws.Range("K137").Value = 64600
$ws.Range("L137").Value = 57100
$ws.Range("M137").Value = -59500
$ws.Range("N137").Value = -67300

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 530000
$ws.Range("J141").Value = 530000
$ws.Range("L141").Value = 530000
$ws.Range("N141").Value = -540360
